$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New topic descriptions that replace the numeric issue ids in column A
$s2 = @'
Update the Local Unit edit workflow
- **Feature:** Allow logged-in RCRC Movement users to edit local unit information.
- **Tasks:**
  - Display a change summary after edits are done and allow the user to "Submit".
  - Update front-end with pending edits validation status after submission.
  - Highlight the fields that have been modified when viewing the form.
'@
$s3 = @'
Update the Local Unit addition workflow
- **Feature:** Allow logged-in RCRC Movement users to add new local units.
- **Key Tasks:**
  - Create a new entry using the Local Unit form
  - Display the new entry with `Not validated` status.
'@
$s4 = @'
Implement the Local Unit deletion workflow
- **Feature:** Allow Local Unit Admin users to delete local units.
- **Key Tasks:**
  - Present deletion options (`Non-existent`, `Incorrectly Added`, `Security Concerns`, `Other`).
  - Require a comment box for additional details.
'@
$s5 = @'
Field Report 3.0
Background
The purpose of the Simplified Reporting Session during the KL meeting was to pow wow around the following fundamental issues. 
1.	Reduce Field/Question Count:
o	Identify and eliminate redundant fields across the field report, and a new Monitoring tool from the Project 3W and Emergency 3W tools.
2.	Standardize Data Interoperability:
o	Ensure field names, data types, and formats are consistent between the platforms.
o	Create a shared field structure to ease the data transfer between tools.
3.	Use Required/Optional Fields:
o	Make fields essential for core reporting mandatory, while others are optional, reducing the burden on the users.
4.	Implement Dynamic Forms:
o	Conditional questions that only appear based on previous answers can streamline the user experience.
Proposed state
![Image](https://github.com/user-attachments/assets/da50b6ef-9ed5-40b8-80e2-f31442aea9e8)

'@

# Drop the confidence_score column (C) entirely
$ws.Columns("C").Delete()

# Replace the issue-id numbers in column A with topic descriptions
$ws.Range("A2").Value = $s2
$ws.Range("A3").Value = $s3
$ws.Range("A4").Value = $s4
$ws.Range("A5").Value = $s5

# The multi-line text triggers Excel auto row-height; restore the default
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()
$ws.Rows(5).AutoFit()

